$d = $word.ActiveDocument
$d.Content.Find.Execute("analise_dados_XX_aaaa", $true, $false, $false, $false, $false,
                         $true, 1, $false, "analise_dados_JF_2021", 2)
